# Daily attendance processing - 2025-12-11 23:28:14
# Swap the order of the first two comma-separated entries in the
# "Recorded By" column (column G) for every data row on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Count -ge 2) {
            $tmp = $parts[0]
            $parts[0] = $parts[1]
            $parts[1] = $tmp
            $cell.Value2 = [string]::Join(", ", $parts)
        }
    }
}
